# Daily roll-forward of the "剩余" (days remaining) tracker sheet.
#
# Each data row (2..99) tracks a supply that lasts D ("总天", total days)
# days, starting on date F ("开始时间"), with E ("剩余") holding how many
# days are left as of "today". This script advances the sheet by one day:
#   - For a normal row, the remaining-day counter E is simply decremented
#     by 1 (one more day has elapsed since the last update).
#   - When a row's counter would reach 0 (i.e. E was already 1, its last
#     day), the supply is treated as replenished instead: E resets back
#     up to the row's total D, and the start date F is set to the new
#     "today" (2025-10-14).
# Row 36 is left untouched because its F value is a malformed/garbled
# date (202510929, 9 digits) that cannot be interpreted as a real date,
# so it is skipped rather than guessed at.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251014

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
if ($lastRow -lt 1) {
    $lastRow = 99
}

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $d = $dCell.Value2
    $e = $eCell.Value2
    $f = $fCell.Value2

    if ($d -eq $null -or $e -eq $null -or $f -eq $null) {
        continue
    }

    $fText = [string]$f
    if ($fText.Length -ne 8) {
        # Malformed start date (e.g. row 36's "202510929") - leave as-is.
        continue
    }

    if ($e -le 1) {
        $eCell.Value = $d
        $fCell.Value = $today
    } else {
        $eCell.Value = $e - 1
    }
}
